$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet2: the size-tracking log. Three new rows of data (97-99) are recorded,
# and the baseline (column D) drops from 21276 to 20036 starting at row 99
# because the previous baseline (21276) was actually taken after some
# changes; this commit records the real/better baseline and carries it
# through the remaining rows (99-153).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

# Row 97: new size data point (no annotation text)
$ws2.Range("B97").Value = 20212

# Row 98: new size data point + "what" / "deltafrombase" annotations
$ws2.Range("A98").Value = "lots of bug fixes, improvements to filters and echo"
$ws2.Range("B98").Value = 20244
$ws2.Range("F98").Value = "impressively small change. There were additions + optimizations and apparently mostly balanced out."

# Row 99: new size data point, and this is where the baseline (column D)
# changes to the corrected value of 20036
$ws2.Range("A99").Value = "this is actually the correct baseline"
$ws2.Range("B99").Value = 20244
$ws2.Range("D99").Value = 20036
$ws2.Range("F99").Value = "so the baseline was after changes; this is a better baseline for the 2023 version. I have added 208 bytes of code. Pretty balanced."

# Rows 100-153: carry the corrected baseline (20036) forward
$ws2.Range("D100:D153").Value = 20036

# ---------------------------------------------------------------------------
# View-state updates (selections / scroll position) to match the saved
# workbook state.
# ---------------------------------------------------------------------------

# Sheet1: selection moves from J15 to J18
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate() | Out-Null
$ws1.Range("J18").Select() | Out-Null

# Sheet2: becomes active again, pane scrolled down one row (A66 -> A67) and
# selection moves from F97 to F100
$ws2.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 67
$win.ScrollColumn = 1
$ws2.Range("F100").Select() | Out-Null
